# Governance Body Suite.xlsx - "GB entity files" update
#
# Summary of the target edit (from the OOXML diff):
#  - "Test Cases" sheet: C3/C4 flip from "N" to "Y", D3 gets a "PASS" result,
#    and the sheet's remembered selection moves from C10 to C5.
#  - "GBCreation" sheet becomes the active tab (was "GBUpdate"); its remembered
#    selection becomes AL2; the sample date fields (X2/Y2) change from 1/30 to
#    2/10, and its Results cell (AI2) flips from "FAIL" to "PASS".
#  - "GBUpdate" sheet is no longer the active tab; its remembered selection
#    moves from G9 to A5.
#  - The now-unused shared string "N" disappears from the shared-strings table
#    (handled automatically by the runtime once no cell references it).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Test Cases" sheet
# ---------------------------------------------------------------------------
$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsTestCases.Range("C3").Value = "Y"
$wsTestCases.Range("D3").Value = "PASS"
$wsTestCases.Range("C4").Value = "Y"

# Update this sheet's saved selection (it is not the active tab in the end,
# so we activate it now and let a later Activate() on another sheet take
# over "tabSelected").
$wsTestCases.Activate()
$wsTestCases.Range("C5").Select()

# ---------------------------------------------------------------------------
# "GBUpdate" sheet
# ---------------------------------------------------------------------------
$wsGBUpdate = $wb.Worksheets.Item("GBUpdate")

# Update this sheet's saved selection; like above, it ends up not being the
# active tab, so it is activated only transiently here.
$wsGBUpdate.Activate()
$wsGBUpdate.Range("A5").Select()

# ---------------------------------------------------------------------------
# "GBCreation" sheet
# ---------------------------------------------------------------------------
$wsGBCreation = $wb.Worksheets.Item("GBCreation")
$wsGBCreation.Range("X2").Value = 2
$wsGBCreation.Range("Y2").Value = 10
$wsGBCreation.Range("AI2").Value = "PASS"

# This is the sheet that ends up active/selected, so activate and select it
# last so "tabSelected"/activeTab stick here.
$wsGBCreation.Activate()
$wsGBCreation.Range("AL2").Select()
